$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = tags
    $val = $cell.Value2

    if ($val -eq "[]") {
        $cell.Value = "aparna"
    }
    elseif ($val -eq "['vip']") {
        $cell.Value = "vip, aparna"
    }
}
